$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the preparation text in F2 with new "Kode Transaksi" value
$ws.Range("F2").Value = "Username : 30711;`nPassword : bni1234;`nRole : 08/09 - Pimpinan Bagian Settlement/Penyelia Settlement;`nKode Transaksi : 998"

# Update KODE_JENIS_TRANSAKSI value (M2) from "090" (text) to 998 (number)
$ws.Range("M2").Value = 998

# Update the selected/visible range in the sheet view
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("G2").Select()
